$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1:L1").EntireColumn.Insert()

$ws.Range("J1").Value = "Flipkart URL"
$ws.Range("K1").Value = "Flipkart Offer Price"
$ws.Range("L1").Value = "Flipkart MRP"

$ws.Range("J1:L1").ColumnWidth = 13.25

$ws.Range("K4").Select() | Out-Null
